# Fruta / hortaliza, semanal
# Insert two new weekly-report rows (148, 149) into the "Pepino ensalada"
# sheet, pushing the existing rows 148:168 down to 150:170.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 148:168 down by two rows to make room for the new entries.
$ws.Rows("148:149").Insert()

# Row 148 - new entry
$ws.Cells.Item(148, 1).Value  = 9
$ws.Cells.Item(148, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(148, 3).Value  = "Metropolitana"
$ws.Cells.Item(148, 4).Value  = 44504
$ws.Cells.Item(148, 5).Value  = 13
$ws.Cells.Item(148, 6).Value  = 100112043
$ws.Cells.Item(148, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(148, 8).Value  = "Sin especificar"
$ws.Cells.Item(148, 9).Value  = "Primera"
$ws.Cells.Item(148, 10).Value = 140
$ws.Cells.Item(148, 11).Value = 9000
$ws.Cells.Item(148, 12).Value = 10000
$ws.Cells.Item(148, 13).Value = 9500
$ws.Cells.Item(148, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(148, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(148, 16).Value = 158
$ws.Cells.Item(148, 17).Value = 60
$ws.Cells.Item(148, 18).Value = "Hortaliza"

# Row 149 - new entry
$ws.Cells.Item(149, 1).Value  = 9
$ws.Cells.Item(149, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(149, 3).Value  = "Metropolitana"
$ws.Cells.Item(149, 4).Value  = 44504
$ws.Cells.Item(149, 5).Value  = 13
$ws.Cells.Item(149, 6).Value  = 100112043
$ws.Cells.Item(149, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(149, 8).Value  = "Sin especificar"
$ws.Cells.Item(149, 9).Value  = "Segunda"
$ws.Cells.Item(149, 10).Value = 79
$ws.Cells.Item(149, 11).Value = 8000
$ws.Cells.Item(149, 12).Value = 8000
$ws.Cells.Item(149, 13).Value = 8000
$ws.Cells.Item(149, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(149, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(149, 16).Value = 80
$ws.Cells.Item(149, 17).Value = 100
$ws.Cells.Item(149, 18).Value = "Hortaliza"
